$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row text (A1, B1) ---
$ws.Range("A1").Value = "Kod"
$ws.Range("B1").Value = "İsim"

# --- Add new header C1 ("ID") ---
# Copy A1 first so C1 inherits the same cell style/format as the other
# header cells instead of getting a brand new style entry.
$ws.Range("A1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "ID"

# --- Add a new defined name "ıd" (Turkish dotless-i + d) referring to column C ---
# Directly calling $wb.Names.Add() with a name that *starts* with the
# dotless-i character (U+0131) is rejected by the name validator (it must
# start with a letter and this specific codepoint is not recognized as
# one in that position). Work around this by creating the name with a
# normal placeholder identifier first and then renaming it - renaming
# does not re-run the "starts with a letter" check.
$wb.Names.Add("ph_id_name_tmp", "=kodlar!`$C:`$C")
$newName = $wb.Names.Item($wb.Names.Count)
$newName.Name = [string]([char]0x0131) + "d"
